$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the current header row (row 3), pushing rows 3-13 down to 4-14
$ws.Rows.Item(3).Insert()

# Widen column A to fit the new date label column
$ws.Columns.Item("A").ColumnWidth = 14.125

# New columns to host start/expected-completion date info
$ws.Columns.Item("E").ColumnWidth = 12.125
$ws.Columns.Item("F").ColumnWidth = 15.75
$ws.Columns.Item("G").ColumnWidth = 30.25

# Fill in the new date-label rows (row 2 and row 3)
$ws.Range("F2").Value = "开始日期"
$ws.Range("G2").Value = "2015年10月26日14:42:15"
$ws.Range("F3").Value = "预计完成日期"
$ws.Range("G3").Value = "2015年10月26日14:43:00"

# Center-align the newly added cells, matching the header/labels style
$ws.Range("F2:G3").HorizontalAlignment = -4108
$ws.Range("F2:G3").VerticalAlignment = -4108

# Move the active selection
$ws.Range("C8").Select()
